# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Applies reconciliation updates: corrects row 182 and appends new
# reconciliation rows 231-265 (New Ndogbong Plateau / Pk8 / Socaver Ndongbong / Total Ndokotti).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($r, $a, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
}

# --- Correct existing row 182 (name + balance figures) ---
$ws.Range("B182").Value = "N A DAMA ADAMA"
$ws.Range("F182").Value = 354745
$ws.Range("G182").Value = 258875.14
$ws.Range("H182").Value = 3.700276604138151

# --- Append new reconciliation rows 231-265 ---
Set-RowData 231 237671646117 'FLORENCE MEDONGUE' 'Rte_5' 'New Ndogbong Plateau' 51881.075 623824 571942.925 12.02411476631893 'Ndogbong'
Set-RowData 232 237674440808 'HOUMI EPSE MVEINGUE GUINDOP TATIANA ETS LE CONTENT' 'Rte_0' 'New Ndogbong Plateau' 5000 6322 1322 1.2644 'Ndogbong'
Set-RowData 233 237675788721 'ETS MOBILE FINANCIAL SERVICES MFS TCHAWE MBOUGA JUDITH FLORE' 'Rte_6' 'New Ndogbong Plateau' 145306.5 569488 424181.5 3.919219030119093 'Ndogbong'
Set-RowData 234 237676036914 'NDENGUE ELOKO SAMUEL HERVE ETS MOBILE FINANCIAL SERVICES MFS' 'Rte_0' 'New Ndogbong Plateau' 29860 92721 62861 3.105190890823845 'Ndogbong'
Set-RowData 235 237679550294 'N A FONATIA' 'Rte_0' 'New Ndogbong Plateau' 14649.6 35722 21072.4 2.438428352992573 'Ndogbong'
Set-RowData 236 237679604574 'KINGUE KOMBI VICTORINE SIDONIE VISION TRADING COMPLEX AND TECHNOLOGIES SARL VISION TRADING COMPLEX' 'Rte_0' 'New Ndogbong Plateau' 10309.09090909091 25643 15333.90909090909 2.487416225749559 'Ndogbong'
Set-RowData 237 237682370358 'JEANNE PRISCA NGO DJON EPSE EBANA ZOE' 'Rte_5' 'New Ndogbong Plateau' 24760 328505 303745 13.26756865912763 'Ndogbong'
Set-RowData 238 237682639044 'JOSEPH KAMGA' 'Rte_7' 'New Ndogbong Plateau' 223294 46373 -176921 0.2076768744346019 'Ndogbong'
Set-RowData 239 237683023087 'FAGHUIE ABIBA' 'Rte_5' 'New Ndogbong Plateau' 119280 373215 253935 3.128898390342052 'Ndogbong'
Set-RowData 240 237650874464 'ASSONFACK VANESSA ATB POINT COM' 'Rte_4' 'Pk8' 189390.9090909091 536391 347000.0909090909 2.832189795036721 'Ndogbong'
Set-RowData 241 237653854849 'TITTI GASTON CLEMENT TOP MOBIL' 'Rte_4' 'Pk8' 18363.75 98719 80355.25 5.375753862909264 'Ndogbong'
Set-RowData 242 237654164073 'NGUIDJOL SIMONE ASTRIDE SPECTRUM SPECTRUM' 'Rte_4' 'Pk8' 30457.27272727273 54178 23720.72727272727 1.778819807181446 'Ndogbong'
Set-RowData 243 237671357520 'CALICE WOTI EPSE DJOMO' 'Rte_0' 'Pk8' 16806.66666666667 15906 -900.6666666666679 0.9464101547005156 'Ndogbong'
Set-RowData 244 237672587687 'Kamaha Tomy Nadine LA NEGRESSE SARL' 'Rte_0' 'Pk8' 43904.9 247714 203809.1 5.642058175738926 'Ndogbong'
Set-RowData 245 237674240552 'AUGUSTINE NGO BAYOI' 'Rte_4' 'Pk8' 20270 56178 35908 2.771484953132708 'Ndogbong'
Set-RowData 246 237675239360 'ERIC MBAH AKEN' 'Rte_4' 'Pk8' 34469.9 107033 72563.10000000001 3.10511489734522 'Ndogbong'
Set-RowData 247 237675396752 'BENEDICTE CHANTAL MANTSANG' 'Rte_4' 'Ndogbong' 34635 121999 87364 3.522419517828786 'Ndogbong'
Set-RowData 248 237675626141 'FLORENCE NGUEFACK' 'Rte_0' 'Pk8' 24500 148545 124045 6.063061224489796 'Ndogbong'
Set-RowData 249 237676840777 'ETP109 ETP' 'Rte_0' 'Pk8' 471586.4666666666 0 -471586.4666666666 0 'Ndogbong'
Set-RowData 250 237677833877 'ISSA ISSYAKOU' 'Rte_4' 'Pk8' 38396.7 561227 522830.3 14.61654256746023 'Ndogbong'
Set-RowData 251 237678854978 'NSAMO NDJOUOHOU MICRANGE ETS MOBILE FINANCIAL SERVICES MFS' 'Rte_5' 'Pk8' 40520 258563 218043 6.381120434353406 'Ndogbong'
Set-RowData 252 237679422591 'ETS LE CONTENT 42' 'Rte_0' 'Pk8' 141149.8 401400 260250.2 2.843787238805865 'Ndogbong'
Set-RowData 253 237650353920 'MENIAPI HELENE EDOSSINE TOP MOBIL TELECOM' 'Rte_4' 'Socaver Ndongbong' 218865 980204 761339 4.478578118931762 'Ndogbong'
Set-RowData 254 237651927448 'charity aben awalah' 'Rte_4' 'Socaver Ndongbong' 128653.3333333333 71305 -57348.33333333331 0.5542413721629185 'Ndogbong'
Set-RowData 255 237653294562 'NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS' 'Rte_5' 'Socaver Ndongbong' 126220 271672 145452 2.152368879733798 'Ndogbong'
Set-RowData 256 237678046498 'MFS SOCAVER' 'Rte_4' 'Socaver Ndongbong' 41943.17 303 -41640.17 0.007224060556224053 'Ndogbong'
Set-RowData 257 237679428698 'ETS LE CONTENT 29' 'Rte_8' 'Socaver Ndongbong' 161427.5 7 -161420.5 0.00004336311966672345 'Ndogbong'
Set-RowData 258 237679551262 'LA NEGRESSE LTDLA CBOX R1 MEGNE JUDITH' 'Rte_5' 'Socaver Ndongbong' 65595 17349 -48246 0.2644866224559799 'Ndogbong'
Set-RowData 259 237680574202 'TOUMEWO SAMUEL' 'Rte_6' 'Socaver Ndongbong' 86620 297874 211254 3.438859385823136 'Ndogbong'
Set-RowData 260 237681118330 'SAHA NDESA JONAS LTDLA_POLAS_OTH_NDOGBONG SERIE' 'Rte_5' 'Socaver Ndongbong' 152199.5714285714 303554 151354.4285714286 1.994447140361762 'Ndogbong'
Set-RowData 261 237674446293 'SYDONIE MAFOMA MESSINE' 'Rte_0' 'Total Ndokotti' 5994.285714285715 10536 4541.714285714285 1.75767397521449 'Ndogbong'
Set-RowData 262 237679085953 'MADELEINE NKOUADJIO' 'Rte_0' 'Total Ndokotti' 8950 24018 15068 2.683575418994413 'Ndogbong'
Set-RowData 263 237681662761 'EMMANUEL EKOLLE ELUMBA' 'Rte_0' 'Total Ndokotti' 29783.82 32326 2542.180000000004 1.085354397118973 'Ndogbong'
Set-RowData 264 237682975726 'SYLVIE-ISABELLE DGANHOU EPSE KOUAHOU' 'Rte_0' 'Total Ndokotti' 29809.16666666666 76155 46345.83333333334 2.554751055324145 'Ndogbong'
Set-RowData 265 237683075075 'ESSOM YOUASSI FRANCK LIONEL STYLE. COM' 'Rte_0' 'Total Ndokotti' 173418.2 4943 -168475.2 0.02850335201264919 'Ndogbong'

Write-Output "Applied row 182 correction and appended rows 231-265."
